# Apply scraped cryptocurrency price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '68.248.34'
$ws.Range('E2').Value = '  +1.11%  '
$ws.Range('D3').Value = '3.347.90'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''582.83'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').Value = '''177.30'
$ws.Range('E6').Value = '  +0.88%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.32%  '
$ws.Range('E9').Value = '  +3.41%  '
$ws.Range('E10').Value = '  +1.04%  '
$ws.Range('D11').Value = '''47.98'
$ws.Range('E11').Value = '  +5.70%  '
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').Value = '''686.97'
$ws.Range('E13').Value = '  +4.26%  '
$ws.Range('D14').Value = '3.889.03'
$ws.Range('E14').Value = '  +0.59%  '
$ws.Range('E15').Value = '  +0.34%  '
$ws.Range('D16').Value = '68.292.17'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('E17').Value = '  +1.16%  '
$ws.Range('D18').Value = '3.360.27'
$ws.Range('E18').Value = '  +1.04%  '
$ws.Range('D19').Value = '''17.42'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('E21').Value = '  +0.59%  '
$ws.Range('D22').Value = '''5.45'
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').Value = '''16.94'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('D24').Value = '''100.17'
$ws.Range('E24').Value = '  +0.84%  '
$ws.Range('D25').Value = '''3.90'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('D27').Value = '''9.50'
$ws.Range('E27').Value = '  +2.86%  '
$ws.Range('D28').Value = '''33.00'
$ws.Range('D29').Value = '''8.50'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('D30').Value = '''6.93'
$ws.Range('E30').Value = '  -6.48%  '
$ws.Range('D31').Value = '''561.75'
$ws.Range('E31').Value = '  -5.02%  '
$ws.Range('D32').Value = '''11.06'
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E33').Value = '  +1.02%  '
$ws.Range('D34').Value = '''57.89'
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '3.709.96'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').Value = '''3.29'
$ws.Range('E37').Value = '  -1.46%  '
$ws.Range('E38').Value = '  +4.27%  '
$ws.Range('D39').Value = '''34.70'
$ws.Range('E39').Value = '  +2.08%  '
$ws.Range('E40').Value = '  +1.59%  '
$ws.Range('E41').Value = '  -0.73%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = '''0.335'
$ws.Range('E42').Value = '  +0.78%  '
$ws.Range('B43').Value = 'PEPE'
$ws.Range('C43').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D43').Value = '0.0₃0670'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('E45').Value = '  +1.04%  '
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').Value = '''131.39'
$ws.Range('E50').Value = '  +3.41%  '
$ws.Range('E51').Value = '  +0.21%  '
